$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M2, N2 and N3 take purely-numeric-looking text ("33", "222", "55"), so
# force a text number format first or Excel would silently coerce them to
# numbers. M3 ("666ert") is already non-numeric and needs no such nudge.
$ws.Range("M2").NumberFormat = "@"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N3").NumberFormat = "@"

$ws.Range("M2").Value = "33"
$ws.Range("N2").Value = "222"
$ws.Range("M3").Value = "666ert"
$ws.Range("N3").Value = "55"
